$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 16 new rows starting at row 13 (pushes old rows 13-18 down to 29-34)
$ws.Rows("13:28").Insert()

# Row 16 / group 4 / A-2827768
$ws.Cells.Item(13,1).Value = 16
$ws.Cells.Item(13,2).Value = 4
$ws.Cells.Item(13,3).Value = "A-2827768"
$ws.Cells.Item(13,4).Value = "0016.png"
$ws.Cells.Item(13,5).Value = 27

$ws.Cells.Item(14,1).Value = 17
$ws.Cells.Item(14,2).Value = 4
$ws.Cells.Item(14,3).Value = "A-2827768"
$ws.Cells.Item(14,4).Value = "0017.png"
$ws.Cells.Item(14,5).Value = 29

$ws.Cells.Item(15,1).Value = 18
$ws.Cells.Item(15,2).Value = 4
$ws.Cells.Item(15,3).Value = "A-2827768"
$ws.Cells.Item(15,4).Value = "0018.png"
$ws.Cells.Item(15,5).Value = 1

$ws.Cells.Item(16,1).Value = 19
$ws.Cells.Item(16,2).Value = 4
$ws.Cells.Item(16,3).Value = "A-2827768"
$ws.Cells.Item(16,4).Value = "0019.png"
$ws.Cells.Item(16,5).Value = 34

# Row 20-22 / group 5 / A-2827812
$ws.Cells.Item(17,1).Value = 20
$ws.Cells.Item(17,2).Value = 5
$ws.Cells.Item(17,3).Value = "A-2827812"
$ws.Cells.Item(17,4).Value = "0020.png"
$ws.Cells.Item(17,5).Value = 33

$ws.Cells.Item(18,1).Value = 21
$ws.Cells.Item(18,2).Value = 5
$ws.Cells.Item(18,3).Value = "A-2827812"
$ws.Cells.Item(18,4).Value = "0021.png"
$ws.Cells.Item(18,5).Value = 24

$ws.Cells.Item(19,1).Value = 22
$ws.Cells.Item(19,2).Value = 5
$ws.Cells.Item(19,3).Value = "A-2827812"
$ws.Cells.Item(19,4).Value = "0022.png"
$ws.Cells.Item(19,5).Value = 24

# Row 23-24 / group 6 / A-2827885
$ws.Cells.Item(20,1).Value = 23
$ws.Cells.Item(20,2).Value = 6
$ws.Cells.Item(20,3).Value = "A-2827885"
$ws.Cells.Item(20,4).Value = "0023.png"
$ws.Cells.Item(20,5).Value = 24

$ws.Cells.Item(21,1).Value = 24
$ws.Cells.Item(21,2).Value = 6
$ws.Cells.Item(21,3).Value = "A-2827885"
$ws.Cells.Item(21,4).Value = "0024.png"
$ws.Cells.Item(21,5).Value = 37

# Row 25-27 / group 7 / A-2828052
$ws.Cells.Item(22,1).Value = 25
$ws.Cells.Item(22,2).Value = 7
$ws.Cells.Item(22,3).Value = "A-2828052"
$ws.Cells.Item(22,4).Value = "0025.png"
$ws.Cells.Item(22,5).Value = 24

$ws.Cells.Item(23,1).Value = 26
$ws.Cells.Item(23,2).Value = 7
$ws.Cells.Item(23,3).Value = "A-2828052"
$ws.Cells.Item(23,4).Value = "0026.png"
$ws.Cells.Item(23,5).Value = 3

$ws.Cells.Item(24,1).Value = 27
$ws.Cells.Item(24,2).Value = 7
$ws.Cells.Item(24,3).Value = "A-2828052"
$ws.Cells.Item(24,4).Value = "0027.png"
$ws.Cells.Item(24,5).Value = 33

# Row 28-30 / group 8 / A-2828231
$ws.Cells.Item(25,1).Value = 28
$ws.Cells.Item(25,2).Value = 8
$ws.Cells.Item(25,3).Value = "A-2828231"
$ws.Cells.Item(25,4).Value = "0028.png"
$ws.Cells.Item(25,5).Value = 36

$ws.Cells.Item(26,1).Value = 29
$ws.Cells.Item(26,2).Value = 8
$ws.Cells.Item(26,3).Value = "A-2828231"
$ws.Cells.Item(26,4).Value = "0029.png"
$ws.Cells.Item(26,5).Value = 14

$ws.Cells.Item(27,1).Value = 30
$ws.Cells.Item(27,2).Value = 8
$ws.Cells.Item(27,3).Value = "A-2828231"
$ws.Cells.Item(27,4).Value = "0030.png"
$ws.Cells.Item(27,5).Value = 37

# Row 31 / group 9 / A-2828314
$ws.Cells.Item(28,1).Value = 31
$ws.Cells.Item(28,2).Value = 9
$ws.Cells.Item(28,3).Value = "A-2828314"
$ws.Cells.Item(28,4).Value = "0031.png"
$ws.Cells.Item(28,5).Value = 29

# Update the view: zoom to 112% and move selection to G26
$excel.ActiveWindow.Zoom = 112
$ws.Range("G26").Select()
